# fix(gui) step 1 and 2
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Step forward the date in A1 by one day (45308 -> 45309)
$ws.Range("A1").Value = 45309

# Update prices in D33:D36
$ws.Range("D33").Value = 445.44
$ws.Range("D34").Value = 487.2
$ws.Range("D35").Value = 546.36
$ws.Range("D36").Value = 664.6799999999999
